$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- View changes ----
$ws.Activate()
$excel.ActiveWindow.Zoom = 110

# ---- H7: "Name" label (bold, centered, top+left thick border) ----
$ws.Range("H7").Value = "Name"
$ws.Range("H7").Font.Bold = $true
$ws.Range("H7").HorizontalAlignment = -4108
$ws.Range("H7").Borders.Item(7).Weight = -4138
$ws.Range("H7").Borders.Item(8).Weight = -4138

# ---- I7: "Rating" label (bold, centered, top+right thick border) ----
$ws.Range("I7").Value = "Rating"
$ws.Range("I7").Font.Bold = $true
$ws.Range("I7").HorizontalAlignment = -4108
$ws.Range("I7").Borders.Item(10).Weight = -4138
$ws.Range("I7").Borders.Item(8).Weight = -4138

# ---- H8: lookup value (bottom+left thick border) ----
$ws.Range("H8").Value = "Groundhog Day"
$ws.Range("H8").Borders.Item(7).Weight = -4138
$ws.Range("H8").Borders.Item(9).Weight = -4138

# ---- I8: VLOOKUP formula (bottom+right thick border) ----
$ws.Range("I8").ClearFormats()
$ws.Range("I8").Formula = "=VLOOKUP(H8,A3:F22,3,FALSE)"
$ws.Range("I8").Borders.Item(10).Weight = -4138
$ws.Range("I8").Borders.Item(9).Weight = -4138

# ---- Column H width (new column used for the VLOOKUP demo) ----
$ws.Columns.Item(8).EntireColumn.AutoFit()

# ---- Row heights for the thick-bottom-border rows ----
$ws.Range("A6").EntireRow.RowHeight = 15.75
$ws.Range("A8").EntireRow.RowHeight = 15.75

# ---- Remove AutoFilter ----
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# ---- Selection ----
$ws.Range("H9").Select()
